$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dumbbell")
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$s1 = $chart.SeriesCollection(1)
try {
  $interior = $s1.Interior
  Write-Host "ThemeColor:" $interior.ThemeColor
  $interior.ThemeColor = 6
  Write-Host "set theme color via Interior"
} catch { Write-Host "ERR:" $_ }
try {
  $border = $s1.Border
  Write-Host "border color:" $border.Color
} catch { Write-Host "ERR2:" $_ }
